# Upload "upah tidak tetap" template: lock the "Kode Jenis Upah" column to
# "Bulanan" (monthly, code 2) only - update the data sample row and the
# instructions sheet text accordingly, and move the active tab / selections
# to match.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsPetunjuk = $wb.Worksheets.Item("Petunjuk")

# Data sheet: "Kode Jenis Upah" sample value changes from 1 (Tahunan) to 2 (Bulanan)
$wsData.Range("D2").Value = "2"

# Petunjuk sheet: instruction text now documents only the monthly (code 2) option
$wsPetunjuk.Range("A6").Value = "5. Kode Jenis Upah = 2. Bulanan"

# Leave the saved selection on the (now inactive) Data sheet
$wsData.Range("C19").Select()

# Petunjuk becomes the active tab, with its own saved selection
$wsPetunjuk.Activate()
$wsPetunjuk.Range("A10").Select()
